$d = $word.ActiveDocument

$d.Content.Find.Execute("literature.ccdc", $false, $false, $false, $false, $false, $true, 1, $false, "literature.ccdc", 2)
$d.Content.Find.Execute("literature.finalcif", $false, $false, $false, $false, $false, $true, 1, $false, "literature.finalcif", 2)
